# Add the new "Pakistan" row (PK / Pakistan / Asia - others (3)) to the SFC
# Country Group lookup table on Sheet1, right after the last existing row
# (Yemen, row 37).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A38").Value = "PK"
$ws.Range("B38").Value = "Pakistan"
$ws.Range("C38").Value = "Asia - others (3)"

# Column C needs to get a little wider to comfortably fit the existing
# group labels now that another row has been added underneath.
$ws.Columns("C").ColumnWidth = 22.45

# Leave the selection on the newly entered cell, matching where the user's
# cursor ended up after typing the last value of the new row.
[void]$ws.Range("C38").Select()
